# Adds the two "Extra" sheets for Afghanistan ODI stats:
#   "ODI Batting Extra"  -> new sheet (sheetId 4)
#   "ODI Bowling Extra"  -> new sheet (sheetId 5)
#
# Both are appended after the existing "ODI Bowling" sheet, reuse the same
# header styling (bold + thin border, centered/top aligned) already used by
# the other sheets' header rows, and store numeric-looking strings
# ("4267", "0", "20.00%") as genuine text, matching the rest of the workbook
# (which stores everything as text except for genuinely numeric fields like
# BATTING_POSITION).

$wb = $excel.ActiveWorkbook

# A scratch cell (far outside any used range) used purely to coerce
# otherwise-numeric-looking strings ("4267", "0", "20.00%") into text before
# copying them (as values only) into the real destination cells. This
# mirrors typing an apostrophe-prefixed value into Excel, without leaving
# any NumberFormat/style behind on the destination cells themselves.
$scratchSheet = $wb.Worksheets.Item(1)
$scratch = $scratchSheet.Cells.Item(500, 500)
$scratch.NumberFormat = "@"

function Set-TextValue($range, $text) {
    $scratch.Value = $text
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Match the page margins already used by the rest of the workbook's sheets
# (0.75in / 0.75in / 1in / 1in / 0.5in / 0.5in), instead of Excel's defaults
# for a freshly inserted sheet.
function Set-StandardMargins($ws) {
    $ws.PageSetup.LeftMargin = 54
    $ws.PageSetup.RightMargin = 54
    $ws.PageSetup.TopMargin = 72
    $ws.PageSetup.BottomMargin = 72
    $ws.PageSetup.HeaderMargin = 36
    $ws.PageSetup.FooterMargin = 36
}

# Used to copy the existing bold/border/centered header style onto the new
# header rows so they reuse the same style record as the rest of the
# workbook instead of minting new font/style entries.
$headerStyleSource = $wb.Worksheets.Item(2).Range("A1")

# ---------------------------------------------------------------------
# Sheet: "ODI Batting Extra"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsBattingExtra = $wb.Worksheets.Add($null, $lastSheet)
$wsBattingExtra.Name = "ODI Batting Extra"
Set-StandardMargins $wsBattingExtra

$headerStyleSource.Copy() | Out-Null
$wsBattingExtra.Range("A1:F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$battingHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $battingHeaders.Count; $i++) {
    Set-TextValue $wsBattingExtra.Cells.Item(1, $i + 1) $battingHeaders[$i]
}

Set-TextValue $wsBattingExtra.Cells.Item(2, 1) "4267"
$wsBattingExtra.Cells.Item(2, 2).Value = 9
Set-TextValue $wsBattingExtra.Cells.Item(2, 3) "0"
Set-TextValue $wsBattingExtra.Cells.Item(2, 4) "0"
Set-TextValue $wsBattingExtra.Cells.Item(2, 6) "NO"

# ---------------------------------------------------------------------
# Sheet: "ODI Bowling Extra"
# ---------------------------------------------------------------------
$wsBowlingExtra = $wb.Worksheets.Add($null, $wsBattingExtra)
$wsBowlingExtra.Name = "ODI Bowling Extra"
Set-StandardMargins $wsBowlingExtra

$headerStyleSource.Copy() | Out-Null
$wsBowlingExtra.Range("A1:C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$bowlingHeaders = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($i = 0; $i -lt $bowlingHeaders.Count; $i++) {
    Set-TextValue $wsBowlingExtra.Cells.Item(1, $i + 1) $bowlingHeaders[$i]
}

Set-TextValue $wsBowlingExtra.Cells.Item(2, 1) "4267"
Set-TextValue $wsBowlingExtra.Cells.Item(2, 2) "0"
Set-TextValue $wsBowlingExtra.Cells.Item(2, 3) "20.00%"

# Clean up the scratch cell so it doesn't show up as used range / content.
$scratch.Clear() | Out-Null

$wsBattingExtra.Range("A1").Select() | Out-Null
